$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---
# Column C (3) narrows to match columns A/B
$ws.Columns.Item(3).ColumnWidth = 2.140625
# Column F (6) widens to match columns D/E
$ws.Columns.Item(6).ColumnWidth = 3.140625
# Columns I, J, K (9-11) now hold decimal/percentage data, widen to match columns M:Q
$ws.Columns.Item(9).ColumnWidth = 5.7109375
$ws.Columns.Item(10).ColumnWidth = 5.7109375
$ws.Columns.Item(11).ColumnWidth = 5.7109375

# --- Row 1 data refresh ---
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 20
$ws.Range("E1").Value = 31
$ws.Range("F1").Value = 22
$ws.Range("G1").Value = 16
$ws.Range("H1").Value = 18
$ws.Range("I1").Value = 0.059
$ws.Range("J1").Value = 0.026
$ws.Range("K1").Value = 0.061
